$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 19.01164980593796
$ws.Range("C2").Value = 5.266203629647237
$ws.Range("D2").Value = 8.475935357025193
$ws.Range("E2").Value = 10.66851330593063
$ws.Range("F2").Value = 42.06301780364806
$ws.Range("K2").Value = 15.21615303703825
$ws.Range("L2").Value = 10.25156721744359
$ws.Range("N2").Value = 23.72442181002087

$ws.Range("B3").Value = 18.82691776739093
$ws.Range("C3").Value = 5.115554408923071
$ws.Range("D3").Value = 8.487262257346675
$ws.Range("E3").Value = 10.67603973725438
$ws.Range("F3").Value = 41.90554509757469
$ws.Range("K3").Value = 15.0915623774294
$ws.Range("L3").Value = 10.2450932523301
$ws.Range("N3").Value = 23.75357729588042

$ws.Range("B4").Value = 18.71794610655805
$ws.Range("C4").Value = 5.019570367264838
$ws.Range("D4").Value = 8.494427539701309
$ws.Range("E4").Value = 10.68232289721531
$ws.Range("F4").Value = 41.81724240539629
$ws.Range("K4").Value = 15.0187262712085
$ws.Range("L4").Value = 10.24307119016055
$ws.Range("N4").Value = 23.7733063517622

$ws.Range("B5").Value = 18.67470715159028
$ws.Range("C5").Value = 4.979606918859745
$ws.Range("D5").Value = 8.497400636014675
$ws.Range("E5").Value = 10.68530135535586
$ws.Range("F5").Value = 41.78338688379536
$ws.Range("K5").Value = 14.98999542398255
$ws.Range("L5").Value = 10.24273923136101
$ws.Range("N5").Value = 23.78180538484959

$ws.Range("B6").Value = 18.66759923966767
$ws.Range("C6").Value = 4.97292057124111
$ws.Range("D6").Value = 8.497897535942657
$ws.Range("E6").Value = 10.68582117493787
$ws.Range("F6").Value = 41.7778942863354
$ws.Range("K6").Value = 14.98528291467452
$ws.Range("L6").Value = 10.24271385093578
$ws.Range("N6").Value = 23.78324437496591

$ws.Range("B7").Value = 18.71735818097938
$ws.Range("C7").Value = 5.019034806837474
$ws.Range("D7").Value = 8.494467420262815
$ws.Range("E7").Value = 10.68236137311085
$ws.Range("F7").Value = 41.81677717554197
$ws.Range("K7").Value = 15.0183349112406
$ws.Range("L7").Value = 10.24306471990628
$ws.Range("N7").Value = 23.77341911347395

$ws.Range("B8").Value = 18.94705903817286
$ws.Range("C8").Value = 5.214998555492086
$ws.Range("D8").Value = 8.479797330255243
$ws.Range("E8").Value = 10.67076359213349
$ws.Range("F8").Value = 42.00699327286502
$ws.Range("K8").Value = 15.17245207988064
$ws.Range("L8").Value = 10.24893030439069
$ws.Range("N8").Value = 23.73409518078321

$ws.Range("B9").Value = 19.43063721122188
$ws.Range("C9").Value = 5.570566680476908
$ws.Range("D9").Value = 8.452687818551837
$ws.Range("E9").Value = 10.66119867400886
$ws.Range("F9").Value = 42.44556865764058
$ws.Range("K9").Value = 15.50236473804646
$ws.Range("L9").Value = 10.27587645258017
$ws.Range("N9").Value = 23.67149296911011

$ws.Range("B10").Value = 19.80303359183175
$ws.Range("C10").Value = 5.813037471092737
$ws.Range("D10").Value = 8.433764101216571
$ws.Range("E10").Value = 10.6621899827886
$ws.Range("F10").Value = 42.80629427370013
$ws.Range("K10").Value = 15.7596832402031
$ws.Range("L10").Value = 10.30500359385082
$ws.Range("N10").Value = 23.63436293133046

$ws.Range("B11").Value = 19.97547845760151
$ws.Range("C11").Value = 5.919047021755796
$ws.Range("D11").Value = 8.425367191686048
$ws.Range("E11").Value = 10.66437671337536
$ws.Range("F11").Value = 42.97839853697953
$ws.Range("K11").Value = 15.8795468328965
$ws.Range("L11").Value = 10.32025629800983
$ws.Range("N11").Value = 23.61939988494012

$ws.Range("B12").Value = 20.04115772976598
$ws.Range("C12").Value = 5.958556980262647
$ws.Range("D12").Value = 8.422217649610317
$ws.Range("E12").Value = 10.66545368373635
$ws.Range("F12").Value = 43.04468764309171
$ws.Range("K12").Value = 15.9253017645944
$ws.Range("L12").Value = 10.32631755224975
$ws.Range("N12").Value = 23.614011253201

$ws.Range("B13").Value = 20.02699664530553
$ws.Range("C13").Value = 5.950076223150455
$ws.Range("D13").Value = 8.422894621163575
$ws.Range("E13").Value = 10.66521068130562
$ws.Range("F13").Value = 43.03036199338288
$ws.Range("K13").Value = 15.91543199915295
$ws.Range("L13").Value = 10.32499950427863
$ws.Range("N13").Value = 23.6151594416405

$ws.Range("B14").Value = 19.98087469046693
$ws.Range("C14").Value = 5.922310314702943
$ws.Range("D14").Value = 8.425107473381271
$ws.Range("E14").Value = 10.66446033242712
$ws.Range("F14").Value = 42.98383000154269
$ws.Range("K14").Value = 15.88330402254117
$ws.Range("L14").Value = 10.32074925730545
$ws.Range("N14").Value = 23.61895099440635

$ws.Range("B15").Value = 19.95267115444552
$ws.Range("C15").Value = 5.905219930628637
$ws.Range("D15").Value = 8.426466833153071
$ws.Range("E15").Value = 10.66403311375265
$ws.Range("F15").Value = 42.95547221205225
$ws.Range("K15").Value = 15.86367111359961
$ws.Range("L15").Value = 10.31818294239301
$ws.Range("N15").Value = 23.62130958573549

$ws.Range("B16").Value = 19.79182004576073
$ws.Range("C16").Value = 5.806021691747886
$ws.Range("D16").Value = 8.434317096071453
$ws.Range("E16").Value = 10.66208194603273
$ws.Range("F16").Value = 42.79520508890722
$ws.Range("K16").Value = 15.7519030900143
$ws.Range("L16").Value = 10.30404687682315
$ws.Range("N16").Value = 23.63537962091695

$ws.Range("B17").Value = 19.69387978337027
$ws.Range("C17").Value = 5.744055304188575
$ws.Range("D17").Value = 8.439186998829653
$ws.Range("E17").Value = 10.66132908270596
$ws.Range("F17").Value = 42.69891353088758
$ws.Range("K17").Value = 15.68402895337469
$ws.Range("L17").Value = 10.29588602616749
$ws.Range("N17").Value = 23.64450509514237

$ws.Range("B18").Value = 19.63783681698181
$ws.Range("C18").Value = 5.708010785781732
$ws.Range("D18").Value = 8.442007967760311
$ws.Range("E18").Value = 10.66105948493474
$ws.Range("F18").Value = 42.6442851645097
$ws.Range("K18").Value = 15.64525621060849
$ws.Range("L18").Value = 10.29138078995753
$ws.Range("N18").Value = 23.64993523949886

$ws.Range("B19").Value = 19.61891317897834
$ws.Range("C19").Value = 5.695738031514306
$ws.Range("D19").Value = 8.44296652962751
$ws.Range("E19").Value = 10.66099629224511
$ws.Range("F19").Value = 42.62591980289431
$ws.Range("K19").Value = 15.6321753939424
$ws.Range("L19").Value = 10.28988787272681
$ws.Range("N19").Value = 23.65180493971762

$ws.Range("B20").Value = 19.70427613742778
$ws.Range("C20").Value = 5.750693575231221
$ws.Range("D20").Value = 8.438666528497514
$ws.Range("E20").Value = 10.66139231668074
$ws.Range("F20").Value = 42.70908593429449
$ws.Range("K20").Value = 15.69122695040263
$ws.Range("L20").Value = 10.29673525397442
$ws.Range("N20").Value = 23.64351489475067

$ws.Range("B21").Value = 19.99441203442766
$ws.Range("C21").Value = 5.930483151286595
$ws.Range("D21").Value = 8.424456687734629
$ws.Range("E21").Value = 10.66467397910365
$ws.Range("F21").Value = 42.99746754738044
$ws.Range("K21").Value = 15.8927311869494
$ws.Range("L21").Value = 10.32198993459262
$ws.Range("N21").Value = 23.61782978797185

$ws.Range("B22").Value = 20.18620938580692
$ws.Range("C22").Value = 6.044286846275622
$ws.Range("D22").Value = 8.415345567364472
$ws.Range("E22").Value = 10.66826912210862
$ws.Range("F22").Value = 43.19243466392798
$ws.Range("K22").Value = 16.02653627494023
$ws.Range("L22").Value = 10.34015737936406
$ws.Range("N22").Value = 23.60266101968242

$ws.Range("B23").Value = 20.08366384722268
$ws.Range("C23").Value = 5.983891097887634
$ws.Range("D23").Value = 8.420192334434249
$ws.Range("E23").Value = 10.66621787855143
$ws.Range("F23").Value = 43.08779487145509
$ws.Range("K23").Value = 15.95494174347678
$ws.Range("L23").Value = 10.33030992193636
$ws.Range("N23").Value = 23.61060871718642

$ws.Range("B24").Value = 19.69957511833579
$ws.Range("C24").Value = 5.747693717433767
$ws.Range("D24").Value = 8.438901767074872
$ws.Range("E24").Value = 10.6613632201089
$ws.Range("F24").Value = 42.70448471127541
$ws.Range("K24").Value = 15.68797195740206
$ws.Range("L24").Value = 10.29635073683319
$ws.Range("N24").Value = 23.64396199185631

$ws.Range("B25").Value = 19.29658223897302
$ws.Range("C25").Value = 5.477589029896144
$ws.Range("D25").Value = 8.459845909300151
$ws.Range("E25").Value = 10.66237652940634
$ws.Range("F25").Value = 42.32004519311511
$ws.Range("K25").Value = 15.41034627501618
$ws.Range("L25").Value = 10.26694110400113
$ws.Range("N25").Value = 23.68687315864315
